$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query text in B2 was edited: the stray blank line right after the
# first MATCH clause was removed, and the trailing "Cohort" column
# (coalesce(co.cohort_description, '') AS `Cohort`) was dropped from the RETURN list.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Male', 'Female'] and demo.neutered_indicator IN ['No']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# Row heights shrink because the CasesTab cell lost the "Cohort" line (and the
# SamplesTab / FilesTab rows were re-measured too).
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 261

# Selection / scroll position moved from C4 to B2.
[void]$ws.Range("B2").Select()
